# Apply the "Product features.xlsx" documentation update:
#  - Expand the "log attempts..." note (F2) with the "when an error (or
#    above) happens" detail, and widen the Notes column (F) to fit it.
#  - Add a new "User Management" story (row 8) describing user CRUD
#    permissions for Administrators.
#  - Update the remembered cell selection on the "product backlog" and
#    "sprint backlog" sheets.

$wb = $excel.ActiveWorkbook

$productBacklog = $wb.Worksheets.Item("product backlog")
$sprintBacklog  = $wb.Worksheets.Item("sprint backlog")

# --- product backlog: widen the Notes column (F) ---------------------------
$productBacklog.Columns.Item(6).ColumnWidth = 45.28515625

# --- product backlog: extend the existing login-logging note (F2) ---------
$productBacklog.Range("F2").Value = "log attempts in database and send via email when an error (or above) happens"

# --- product backlog: new row 8 - user management story -------------------
$productBacklog.Range("B8").Value = "User Management"
$productBacklog.Range("C8").Value = "Administrator"
$productBacklog.Range("D8").Value = "be able to create, delete and edit other users. The simple users will only be able to edit their data (except from their category)"
$productBacklog.Rows.Item(8).RowHeight = 60

# --- selections --------------------------------------------------------
# Select sprint backlog first so that "product backlog" ends up as the
# active sheet (tabSelected) afterwards, matching the saved workbook state.
$sprintBacklog.Range("E3").Select() | Out-Null
$productBacklog.Range("F3").Select() | Out-Null
